$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '90.842.07'
$ws.Range("E2").Value = '  +1.24%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.154.75'
$ws.Range("E3").Value = '  +3.48%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.88'
$ws.Range("E5").Value = '  +1.53%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '626.62'
$ws.Range("E6").Value = '  +2.27%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.13'
$ws.Range("E7").Value = '  +28.24%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.371'
$ws.Range("E8").Value = '  +2.42%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.00'
$ws.Range("E9").Value = '  -0.05%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.150.20'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.754'
$ws.Range("E11").Value = '  +11.22%  '
$ws.Range("E12").Value = '  +8.33%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.69'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000245'
$ws.Range("E14").Value = '  +2.32%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '34.92'
$ws.Range("E15").Value = '  +8.67%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '90.547.96'
$ws.Range("E16").Value = '  +1.08%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.731.34'
$ws.Range("E17").Value = '  +3.91%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.142.46'
$ws.Range("E18").Value = '  +3.18%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.69'
$ws.Range("E19").Value = '  +11.08%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.28'
$ws.Range("E20").Value = '  +6.47%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '465.61'
$ws.Range("E21").Value = '  +9.41%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0000211'
$ws.Range("E22").Value = '  -2.67%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.09'
$ws.Range("E23").Value = '  +10.56%  '
$ws.Range("E24").Value = '  +5.07%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.88'
$ws.Range("E25").Value = '  +9.26%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '93.08'
$ws.Range("E26").Value = '  +10.94%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.11'
$ws.Range("E27").Value = '  +3.78%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.315.03'
$ws.Range("E28").Value = '  +3.57%  '
$ws.Range("E29").Value = '  -0.01%  '
$ws.Range("E30").Value = '  +2.89%  '
$ws.Range("E31").Value = '  -0.59%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '9.13'
$ws.Range("E32").Value = '  +11.22%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '27.25'
$ws.Range("E33").Value = '  +19.14%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '517.09'
$ws.Range("E34").Value = '  +2.38%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.180'
$ws.Range("E35").Value = '  +31.56%  '
$ws.Range("E36").Value = '  +7.14%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.63'
$ws.Range("E37").Value = '  -2.24%  '
$ws.Range("B38").Value = 'Kaspa'
$ws.Range("C38").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.141'
$ws.Range("E38").Value = '  +7.21%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.88'
$ws.Range("E39").Value = '  +3.76%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.30'
$ws.Range("E40").Value = '  +5.43%  '
$ws.Range("B41").Value = 'Hedera'
$ws.Range("C41").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0866'
$ws.Range("E41").Value = '  +25.86%  '
$ws.Range("B42").Value = 'WhiteBITCoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '22.21'
$ws.Range("E42").Value = '  -0.22%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.415'
$ws.Range("E44").Value = '  +14.71%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.97'
$ws.Range("E45").Value = '  +7.70%  '
$ws.Range("E46").Value = '  +0.00%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '149.86'
$ws.Range("E47").Value = '  +1.89%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '45.62'
$ws.Range("E48").Value = '  +5.38%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.54'
$ws.Range("E49").Value = '  +11.03%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.36'
$ws.Range("E50").Value = '  +11.84%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.675'
$ws.Range("E51").Value = '  +15.00%  '
